$wb = $excel.ActiveWorkbook
$wsSteps = $wb.Worksheets.Item("TestSteps")
$wsData = $wb.Worksheets.Item("TestData")

# Revert row 9 on TestSteps back to jsClick / getData=SetStatus
$wsSteps.Range("A9").Value = "jsClick"
$wsSteps.Range("C9").Value = "'getData=SetStatus"

# Remove the "Radio" column (J) that was added to TestData
$wsData.Range("J1:J2").Delete()

# Revert the year value back to 2024
$wsData.Range("E2").Value = 2024

# Restore original selection / active sheet state
$wsData.Range("F10").Select() | Out-Null
$wsSteps.Activate() | Out-Null
$wsSteps.Range("B3").Select() | Out-Null
